$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits right after "Thief" (before the ", Bug"
# run) in the Helpers paragraph. It needs to move to the end of the following
# "Defense types" paragraph, right after "Nort, " (once "pop up blocker" is
# removed from there).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the (first) "Defense types" paragraph and, within it, the "Nort, "
# text that immediately precedes "pop up blocker".
$f = $d.Content
$f.Find.ClearFormatting()
$f.Find.Execute("Nort, pop up blocker", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Re-home the bookmark right after "Nort, " *before* deleting "pop up blocker"
# so it stays anchored to the surviving text once the deletion happens.
$afterComma = $d.Range($f.Start + 6, $f.Start + 6)   # "Nort, " is 6 characters
$d.Bookmarks.Add("_GoBack", $afterComma)

# Now remove "pop up blocker", leaving "Nort, " behind (only within this one
# paragraph, found via the same match range).
$popUpBlocker = $d.Range($f.Start + 6, $f.End)
$popUpBlocker.Text = ""
